# v3.0 update FCI 27/1/2023
# Rebuild the sheet with the new data order (alphabetical fund rows, with
# "avg" and "total" moved to the bottom) and add a new column C holding the
# 13-01-2023 snapshot next to the existing 06-01-2023 column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row order (label, B-value, C-value), rows 2..19.
$data = @(
    @("Adcap IOL Acciones Argentina", 2771.33,  3769.95),
    @("Bull Market",                  9397.280000000001, 9965.25),
    @("CMA acciones",                 14070.6,  14055.35),
    @("Delta Acciones",               24500,    24468.72),
    @("Delta Internacional",          359.58,   359.72),
    @("Delta Latinoamerica",          974.99,   974.83),
    @("Delta Recursos Naturales",     71663.11, 71627.84),
    @("FBA Acciones Argentinas",      10260.35, 6666.05),
    @("FBA Calificado",               10142.32, 6728.6),
    @("Fima PB Acciones",             114.67,   114.8),
    @("Gainvest Renta Variable",      41083.81, 41102.1),
    @("Goal acciones plus",           3501.03,  3499.52),
    @("Lombardi",                     704.72,   695.95),
    @("MAF",                          10125.97, 10122.59),
    @("Pionero Acciones",             4052.68,  11639.92),
    @("Superfondo ",                  14407.88, 33192.71),
    @("avg",                          13633.15, 14936.49),
    @("total",                        218130.32, 238983.9)
)

function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

function Set-DateTextCell($cell, $text) {
    # Force text interpretation so date-shaped strings like "06-01-2023"
    # are not auto-converted into date serials by the smart-entry parser.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 1: the two date headers, in B1 / C1.
$b1 = $ws.Cells.Item(1, 2)
Set-DateTextCell $b1 "06-01-2023"
Set-HeaderStyle $b1

$c1 = $ws.Cells.Item(1, 3)
Set-DateTextCell $c1 "13-01-2023"
Set-HeaderStyle $c1

# Rows 2..19: label in column A, values in B and C.
$r = 2
foreach ($row in $data) {
    $label = $row[0]
    $bVal = $row[1]
    $cVal = $row[2]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $label
    Set-HeaderStyle $aCell

    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal

    $r++
}
